# EmployeeDetail.xlsx edit:
#  - Rename the "department" header (Sheet1!M1) to "deptId"
#  - Convert the department column (Sheet1!M2:M11) from the text label
#    "Comp.Engg." to a numeric department id (4 for rows 2-6, 7 for rows 7-11)
#  - Leave the selection on M11 (last edited cell) like the authoring session did

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header rename: department -> deptId
$ws.Range("M1").Value = "deptId"

# Replace the textual department label with numeric department ids
$deptIds = @(4, 4, 4, 4, 4, 7, 7, 7, 7, 7)
for ($i = 0; $i -lt $deptIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 13).Value = $deptIds[$i]
}

# Match the final selection/scroll position left by the edit
$ws.Activate()
$ws.Range("M11").Select()

Write-Output "done"
